# feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" worksheet (fund-holding detail, same layout as the
# other quarterly sheets) right before the "总计" (totals) sheet, then
# refresh "总计" with a new leading row summarising the 2022-Q1 quarter
# (existing rows shift down by one).

# Helper: write a value that must stay TEXT even though it looks like a
# number (fund codes such as "011685", position figures like "0.73").
# Plain `.Value = "0.73"` gets auto-coerced to a number by Excel, so we
# round-trip it through a throw-away formula cell + paste-values, which
# lands a genuine string with no extra number-format/style baggage.
function Set-TextValue($ws, $row, $col, $text) {
    $scratch = $ws.Cells.Item(1000, 1)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted just before "总计".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Match the header/index-column look of the other quarterly sheets
# (bold, centred, thin-bordered header; same style on the row-index
# column) by copying their formats over.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$template.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

$q1.Cells.Item(2, 1).Value = 0
Set-TextValue $q1 2 2 "011685"
Set-TextValue $q1 2 3 "创金合信先进装备股票A"
Set-TextValue $q1 2 4 "0.73"
Set-TextValue $q1 2 5 "92.01"
Set-TextValue $q1 2 6 "9.02"
Set-TextValue $q1 2 7 "0.0658"
$q1.Cells.Item(2, 8).Value = 4

$q1.Cells.Item(3, 1).Value = 1
Set-TextValue $q1 3 2 "011686"
Set-TextValue $q1 3 3 "创金合信先进装备股票C"
Set-TextValue $q1 3 4 "0.17"
Set-TextValue $q1 3 5 "92.01"
Set-TextValue $q1 3 6 "9.02"
Set-TextValue $q1 3 7 "0.0153"
$q1.Cells.Item(3, 8).Value = 4

# ---------------------------------------------------------------------
# 2. Refresh "总计" — add a new leading row for 2022-Q1, pushing the
#    existing quarters down one row. Re-fetch the sheet by name: adding
#    "2022-Q1" shifted "总计" along the tab order, so the earlier
#    $totalSheet handle is now stale.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Give new row 7 (A7) the same style as the existing index column before
# the cascade below so it matches A2:A6.
$totalSheet.Cells.Item(6, 1).Copy()
$totalSheet.Cells.Item(7, 1).PasteSpecial(-4122)  # xlPasteFormats

$existingRows = @(
    @("2021-Q4", 5, 0.79),
    @("2021-Q3", 8, 3.12),
    @("2021-Q2", 7, 2.72),
    @("2021-Q1", 6, 4.94),
    @("2020-Q4", 4, 1.34)
)

for ($i = $existingRows.Length - 1; $i -ge 0; $i--) {
    $row = $i + 3
    $totalSheet.Cells.Item($row, 2).Value = $existingRows[$i][0]
    $totalSheet.Cells.Item($row, 3).Value = $existingRows[$i][1]
    $totalSheet.Cells.Item($row, 4).Value = $existingRows[$i][2]
}

for ($row = 2; $row -le 7; $row++) {
    $totalSheet.Cells.Item($row, 1).Value = $row - 2
}

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.08
